# Updates the training schedule worksheet:
#  - Normalizes the x_nrSteps column (F) values from "-0" to "0" for all data rows.
#  - For rows where alienID (H) was 14, corrects the trial parameters:
#      y_corrSteps (E) decreases by 1, y_nrSteps (G) becomes -3, alienID (H) becomes 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    # Column F (x_nrSteps): clear the negative-zero formatting artifact.
    $ws.Cells.Item($r, 6).Value = 0

    # Column H (alienID): rows that were 14 get corrected to 13, with matching
    # adjustments to columns E (y_corrSteps) and G (y_nrSteps).
    if ($ws.Cells.Item($r, 8).Value2 -eq 14) {
        $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 5).Value2 - 1
        $ws.Cells.Item($r, 7).Value = -3
        $ws.Cells.Item($r, 8).Value = 13
    }
}

# Leave the selection on row 15, matching where the edits were last reviewed.
$null = $ws.Range("A15").Select()
